$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "New Billing" label to "New_Billing"
$ws.Range("A2").Value = "New_Billing"

# Move the selection to A3
$ws.Range("A3").Select()
